# Update gh-pages to output generated at 456a3b4
# Refresh "want to go" counts (column F) across the 展览 (Exhibition),
# 本地生活 (Local Life) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    4  = 621
    7  = 404
    8  = 710
    10 = 6966
    12 = 5060
    14 = 6599
    15 = 8198
    20 = 627
    21 = 110
    26 = 1133
    28 = 1563
    29 = 640
    30 = 788
    36 = 1354
    40 = 2831
    42 = 177
    43 = 24
    45 = 482
    47 = 681
    48 = 133
    49 = 4018
}
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# --- Sheet "本地生活" (Local Life) ---
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Cells.Item(2, 6).Value = 4838

# --- Sheet "全部类型" (All Types) ---
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    2  = 4838
    14 = 404
    15 = 710
    19 = 5060
    20 = 6599
    21 = 6599
    27 = 627
    28 = 110
    30 = 1133
    31 = 1563
    32 = 640
    33 = 788
    46 = 24
    48 = 681
    49 = 133
    51 = 4018
}
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
